# Fix bug in kraken2 scores table: insert a new "DADA2Tax" group of 3 rows
# (50%, 30%, 70%) before the existing "Kraken_0.0" group on the CO1 sheet,
# shifting the Kraken_0.0 / Kraken_0.05 / Kraken_0.1 / ... / VSEARCH groups
# down by 3 rows (old rows 98-124 become new rows 101-127).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows starting at row 98 (pushes current row 98 and below down by 3).
$ws.Range("A98:A100").EntireRow.Insert()

# Populate the 3 new rows with the DADA2Tax data.
$data = @(
    @("CO1", "DADA2Tax", "50%", 12.4, 21.5, 0, 65.09999999999999, 0.16, 0.37, 0.22, 0.29, 0.13),
    @("CO1", "DADA2Tax", "30%", 15,   27.4, 0, 56.6,               0.21, 0.35, 0.26, 0.31, 0.15),
    @("CO1", "DADA2Tax", "70%", 12.3, 7.1,  0, 79.59999999999999, 0.13, 0.63, 0.22, 0.36, 0.12)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = 98 + $i
    $row = $data[$i]
    # Column C holds text like "50%" - force text format so Excel doesn't
    # auto-convert it to a numeric percentage value.
    $ws.Cells.Item($r, 3).NumberFormat = "@"
    for ($c = 1; $c -le $row.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
